$d = $word.ActiveDocument

# 1. Change program code "C24B05" -> "PC05"
$d.Content.Find.Execute("C24B05", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "PC05", 2)

# 2. Merge the "The content generated by AI tools are not retrievable..." run
#    (remove the proofErr-induced run split: "The content generated by AI tools " + "are" + " not retrievable...")
$d.Content.Find.Execute("The content generated by AI tools are not retrievable", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The content generated by AI tools are not retrievable", 2)

# 3. Merge "Streamlit" + " link" into a single run "Streamlit link" and
#    drop the now-stale spellStart/spellEnd proofErr markers that wrapped
#    "Streamlit". A plain Find/Replace cannot clear the leading spellStart
#    marker (it sits just before the matched range), so rebuild the whole
#    paragraph via InsertXML, which guarantees the exact target markup.
$rng3 = $d.Content
$rng3.Find.Execute("Streamlit link", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
if ($rng3.Find.Found) {
  $streamlitXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="29CB5239" w14:textId="77777777" w:rsidR="00B67A4B" w:rsidRPr="004E5BCF" w:rsidRDefault="00B67A4B" w:rsidP="000968D3">
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
              </w:rPr>
              <w:t>Streamlit link</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
  $rng3.InsertXML($streamlitXml)
}
